$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" values formatted as plain text in the source data
# (e.g. "67.953.72", "0.636"). Force text format before assignment so purely
# numeric-looking strings are not auto-coerced into floating point numbers,
# then restore the default (unstyled) cell style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.445.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.725.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.882.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.636"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.93%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.721"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000304"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.303.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.717.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.93%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.25%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.127"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.155.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "413.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.63%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "45.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.31%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.19%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0951"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.22%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "66.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "593.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.409"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.53%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.18%  "

$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("E43").Value = "  -6.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0447"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -13.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.53%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.24%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000273"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.20%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.755.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.98%  "
